{"js": "// Insert three new bullet paragraphs describing Siege Analytics achievements\n// (voter file discovery / boundary estimation / cost savings metrics)\n// immediately after the \"Product Development and Platform Architecture\"\n// paragraph in the PARTNER - Siege Analytics section, before the existing\n// bullet list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Product Development and Platform Architecture\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find anchor paragraph \"' + anchorText + '\"');\n}\n\nconst newBullets = [\n  \"\u2022 Conceived and architected redistricting platform incorporating boundary estimation algorithm used by 2,500+ analysts\",\n  \"\u2022 Built multi-tenant data warehouse tracking decades of demographic data, enabling discovery of 500,000+ mischaracterized voters\",\n  \"\u2022 Platform democratized redistricting analysis, reducing costs by 75% and enabling 200+ smaller organizations to participate\"\n];\n\n// Insert after the anchor, one at a time, always right after the anchor so\n// that the final order matches the array order (each new paragraph becomes\n// the new insertion point for the next one).\nlet insertAfter = anchor;\nfor (const text of newBullets) {\n  insertAfter = insertAfter.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs describing Siege Analytics achievements\n# (voter file discovery / boundary estimation / cost savings metrics)\n# immediately after the \"Product Development and Platform Architecture\"\n# paragraph in the PARTNER - Siege Analytics section, before the existing\n# bullet list.\n\n$d = $word.ActiveDocument\n\n$bullets = @(\n    \"\u2022 Conceived and architected redistricting platform incorporating boundary estimation algorithm used by 2,500+ analysts\",\n    \"\u2022 Built multi-tenant data warehouse tracking decades of demographic data, enabling discovery of 500,000+ mischaracterized voters\",\n    \"\u2022 Platform democratized redistricting analysis, reducing costs by 75% and enabling 200+ smaller organizations to participate\"\n)\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Product Development and Platform Architecture*\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find anchor paragraph 'Product Development and Platform Architecture'\"\n}\n\n# Insert each new bullet right after the anchor, one at a time, advancing the\n# insertion point so the bullets end up in the same order as the array.\n$current = $anchor\nforeach ($bulletText in $bullets) {\n    $current.Range.InsertParagraphAfter()\n    $current = $current.Next()\n    $current.Range.Text = $bulletText\n}\n"}
